# feat: add 2022-Q1 data
#
# - Adds a new "2022-Q1" fund-holding sheet (same layout as the other
#   quarterly sheets), placed right before the "总计" summary sheet.
# - Inserts a new leading row into "总计" for the 2022-Q1 totals and
#   shifts the existing history rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0. Drop the existing "总计" sheet; it is rebuilt at the end so that the
#    newly-added sheets keep the same tab order / id sequencing that a
#    plain "append a sheet, then append another" session would produce.
# ---------------------------------------------------------------------
$oldTotals = $wb.Worksheets.Item("总计")
$oldTotals.Delete()

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4" (same
#    layout/headers; only the data row needs a few values changed),
#    appended right after "2021-Q4".
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$refSheet.Copy($null, $refSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Update only the cells that differ from the "2021-Q4" source (the
# header row and A2/B2/C2 already carry the right values via the copy).
$textRange = $newSheet.Range("D2:G2")
$textRange.NumberFormat = "@"
$newSheet.Range("D2").Value = "1.35"
$newSheet.Range("E2").Value = "88.71"
$newSheet.Range("F2").Value = "4.55"
$newSheet.Range("G2").Value = "0.0614"
$textRange.Style = "Normal"

$newSheet.Range("H2").Value = 2

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Rebuild the "总计" (totals) sheet after "2022-Q1": duplicate a
#    quarterly sheet (to inherit its sheet-level setup), wipe its
#    contents, then write the totals table with a new leading row for
#    "2022-Q1" followed by the shifted-down history.
# ---------------------------------------------------------------------
$newSheet.Copy($null, $newSheet)
$totals = $wb.Worksheets.Item("2022-Q1 (2)")
$totals.Name = "总计"
$totals.Range("A1:H10").Clear()

$totals.Range("B1").Value = "日期"
$totals.Range("C1").Value = "持有数量(只)"
$totals.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @{ A = 0; B = "2022-Q1"; C = 1; D = 0.06 },
    @{ A = 1; B = "2021-Q4"; C = 1; D = 0.05 },
    @{ A = 2; B = "2021-Q3"; C = 1; D = 0.04 },
    @{ A = 3; B = "2021-Q2"; C = 1; D = 0.05 },
    @{ A = 4; B = "2021-Q1"; C = 1; D = 0.06 },
    @{ A = 5; B = "2020-Q4"; C = 1; D = 0.06 }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $totals.Range("A$r").Value = $row.A
    $totals.Range("B$r").Value = $row.B
    $totals.Range("C$r").Value = $row.C
    $totals.Range("D$r").Value = $row.D
}

# Mirror the header/index styling (bold + bordered, centered) from the
# quarterly sheets onto the totals header row and the "A" index column.
$refSheet.Range("B1:D1").Copy($totals.Range("B1"))
$refSheet.Range("A2").Copy($totals.Range("A2"))
$refSheet.Range("A2").Copy($totals.Range("A3"))
$refSheet.Range("A2").Copy($totals.Range("A4"))
$refSheet.Range("A2").Copy($totals.Range("A5"))
$refSheet.Range("A2").Copy($totals.Range("A6"))
$refSheet.Range("A2").Copy($totals.Range("A7"))

$totals.Range("A1").Select()

# ---------------------------------------------------------------------
# Restore the originally active sheet/selection.
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$firstSheet.Select()
$firstSheet.Range("A1").Select()
